# Applies the "Added Dashboard of all best insights" edit:
#   1. "1.Data Cleaning.ipynb" run split (proofErr spell markers around "Cleaning.ipynb")
#   2. "app.py" paragraph rewritten to the Dashboard_Best blurb, split into
#      proofErr-wrapped runs the way Word's proofing pass would leave them
#   3. "...except "1.Data Cleaning"" run split (proofErr markers around
#      "Jupyter" and "1.Data")
#
# Runs are rebuilt via Range.InsertXML so the w:proofErr siblings (which
# aren't reachable through the plain Word object-model properties) land
# exactly where Word's spell/grammar checker would place them.
#
# Range.InsertXML deletes the target range and re-inserts the new content;
# when other run content in the same paragraph trails *after* the target
# range, the new content is observed to land at the end of the paragraph
# instead of in place. To avoid that, every InsertXML call below targets a
# range that runs from its start point through the paragraph end (just
# before the paragraph mark), carrying any unchanged trailing text along
# as plain runs, so the insertion point is never ambiguous.

$d = $word.ActiveDocument

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'

function Run-Text($text, $preserve) {
    if ($preserve) {
        return '<w:r>' + $rPr + '<w:t xml:space="preserve">' + $text + '</w:t></w:r>'
    } else {
        return '<w:r>' + $rPr + '<w:t>' + $text + '</w:t></w:r>'
    }
}

function Insert-RunsFrom($startPos, $endPos, $runsXml) {
    $r = $d.Range($startPos, $endPos)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runsXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

$quoteOpen = [char]8220
$quoteClose = [char]8221
$enDash = [char]8211

# --- Paragraph 2: "1.Data Cleaning.ipynb" + " - Category 1 queries" ---
# Split "1.Data Cleaning.ipynb" into "1.Data " + proofErr(spellStart/End)"Cleaning.ipynb",
# keep the trailing " - Category 1 queries" run as-is (carried along verbatim
# since it trails the edited span within the same paragraph).
$full2 = $d.Paragraphs.Item(2).Range
$start2 = $full2.Start
$end2 = $full2.End - 1
$runs2 = (Run-Text "1.Data " $true) +
         '<w:proofErr w:type="spellStart"/>' +
         (Run-Text "Cleaning.ipynb" $false) +
         '<w:proofErr w:type="spellEnd"/>' +
         (Run-Text (" " + $enDash + " Category 1 queries") $true)
Insert-RunsFrom $start2 $end2 $runs2

# --- Paragraph 5: "app.py" + " - This is the streamlit dashboard..." ---
$full5 = $d.Paragraphs.Item(5).Range
$start5 = $full5.Start
$end5 = $full5.End - 1
$runs5 = '<w:proofErr w:type="spellStart"/>' +
         (Run-Text "Dashboard_Best" $false) +
         '<w:proofErr w:type="spellEnd"/>' +
         (Run-Text (" " + $enDash + " This is the ") $true) +
         '<w:proofErr w:type="spellStart"/>' +
         (Run-Text "streamlit" $false) +
         '<w:proofErr w:type="spellEnd"/>' +
         (Run-Text " dashboard. Please run from local " $true) +
         '<w:proofErr w:type="spellStart"/>' +
         (Run-Text "cmd" $false) +
         '<w:proofErr w:type="spellEnd"/>' +
         (Run-Text " prompt" $true)
Insert-RunsFrom $start5 $end5 $runs5

# --- Paragraph 9: "...Final file used in all other Jupyter Notebooks except "1.Data Cleaning"" ---
# Leave the leading "GDM_Python_Aug2025.xlsx" run untouched; rebuild only
# the descriptive run that follows it, through to the paragraph end.
$full9 = $d.Paragraphs.Item(9).Range
$prefixLen = "GDM_Python_Aug2025.xlsx".Length
$start9 = $full9.Start + $prefixLen
$end9 = $full9.End - 1
$runs9 = (Run-Text (" " + $enDash + " Final file used in all other ") $true) +
         '<w:proofErr w:type="spellStart"/>' +
         (Run-Text "Jupyter" $false) +
         '<w:proofErr w:type="spellEnd"/>' +
         (Run-Text (" Notebooks except " + $quoteOpen) $true) +
         '<w:proofErr w:type="gramStart"/>' +
         (Run-Text "1.Data" $false) +
         '<w:proofErr w:type="gramEnd"/>' +
         (Run-Text (" Cleaning" + $quoteClose) $true)
Insert-RunsFrom $start9 $end9 $runs9

Write-Output "DONE"
